# Updates vm_pu results for the 380 kV case: slack bus setpoint
# lowered from 1.05 pu to 1.02 pu, with corresponding bus voltage
# magnitude results recomputed for rows 2-25 (columns B-F, I-M).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.04066742710952
$ws.Range("D2").Value = 1.042112379478265
$ws.Range("E2").Value = 1.038958755264215
$ws.Range("F2").Value = 1.039424397813518
$ws.Range("I2").Value = 1.041318843429398
$ws.Range("J2").Value = 1.045752794676574
$ws.Range("K2").Value = 1.044889987376644
$ws.Range("L2").Value = 1.041745306298794
$ws.Range("M2").Value = 1.042209624658859
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.041984716831792
$ws.Range("D3").Value = 1.043114880628448
$ws.Range("E3").Value = 1.040095706977138
$ws.Range("F3").Value = 1.041371098644693
$ws.Range("I3").Value = 1.04177066314492
$ws.Range("J3").Value = 1.046713961725309
$ws.Range("K3").Value = 1.045702751054996
$ws.Range("L3").Value = 1.042691515913576
$ws.Range("M3").Value = 1.043963548044885
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.042835418703602
$ws.Range("D4").Value = 1.043762072289576
$ws.Range("E4").Value = 1.040830049940091
$ws.Range("F4").Value = 1.042628803438121
$ws.Range("I4").Value = 1.042060812318833
$ws.Range("J4").Value = 1.047333774353136
$ws.Range("K4").Value = 1.046226564695276
$ws.Range("L4").Value = 1.04330186099141
$ws.Range("M4").Value = 1.045096119351155
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.043192660025108
$ws.Range("D5").Value = 1.044033798447362
$ws.Range("E5").Value = 1.041138451802653
$ws.Range("F5").Value = 1.043157092659726
$ws.Range("I5").Value = 1.042182265077159
$ws.Range("J5").Value = 1.047593839600363
$ws.Range("K5").Value = 1.04644627771751
$ws.Range("L5").Value = 1.043557996199931
$ws.Range("M5").Value = 1.045571706786447
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.043252619519142
$ws.Range("D6").Value = 1.044079401878085
$ws.Range("E6").Value = 1.041190215453554
$ws.Range("F6").Value = 1.043245768790876
$ws.Range("I6").Value = 1.042202626766662
$ws.Range("J6").Value = 1.047637476295207
$ws.Range("K6").Value = 1.046483139406331
$ws.Range("L6").Value = 1.043600975952236
$ws.Range("M6").Value = 1.045651528403995
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.042840193718517
$ws.Range("D7").Value = 1.043765704490335
$ws.Range("E7").Value = 1.040834172055964
$ws.Range("F7").Value = 1.042635864214067
$ws.Range("I7").Value = 1.042062437239814
$ws.Range("J7").Value = 1.047337251331654
$ws.Range("K7").Value = 1.04622950246395
$ws.Range("L7").Value = 1.043305285260078
$ws.Range("M7").Value = 1.045102476300954
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.041112961717129
$ws.Range("D8").Value = 1.042451491333662
$ws.Range("E8").Value = 1.039343274111386
$ws.Range("F8").Value = 1.040082705246498
$ws.Range("I8").Value = 1.041471996936039
$ws.Range("J8").Value = 1.046078069050886
$ws.Range("K8").Value = 1.045165101793125
$ws.Range("L8").Value = 1.042065481842463
$ws.Range("M8").Value = 1.042802862676393
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.038056246265958
$ws.Range("D9").Value = 1.040124043781676
$ws.Range("E9").Value = 1.036705623539755
$ws.Range("F9").Value = 1.035568184804161
$ws.Range("I9").Value = 1.04041452839692
$ws.Range("J9").Value = 1.043842708239259
$ws.Range("K9").Value = 1.043273226828707
$ws.Range("L9").Value = 1.039865901199612
$ws.Range("M9").Value = 1.038732171862284
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.036009196029882
$ws.Range("D10").Value = 1.038564310620094
$ws.Range("E10").Value = 1.034939808845737
$ws.Range("F10").Value = 1.032547133898604
$ws.Range("I10").Value = 1.039697919846802
$ws.Range("J10").Value = 1.042341032273726
$ws.Range("K10").Value = 1.042000772164471
$ws.Range("L10").Value = 1.038389199639642
$ws.Range("M10").Value = 1.036005114688604
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.035120513878521
$ws.Range("D11").Value = 1.037886944626095
$ws.Range("E11").Value = 1.034173368244941
$ws.Range("F11").Value = 1.031236080923097
$ws.Range("I11").Value = 1.039384822669035
$ws.Range("J11").Value = 1.041688003951534
$ws.Range("K11").Value = 1.041447067030393
$ws.Range("L11").Value = 1.037747256663494
$ws.Range("M11").Value = 1.034820935432875
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.034790064527949
$ws.Range("D12").Value = 1.03763503606742
$ws.Range("E12").Value = 1.033888396607717
$ws.Range("F12").Value = 1.030748641132396
$ws.Range("I12").Value = 1.03926810015666
$ws.Range("J12").Value = 1.041445014134922
$ws.Range("K12").Value = 1.041240981744086
$ws.Range("L12").Value = 1.037508425649845
$ws.Range("M12").Value = 1.034380559685232
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.034860963153894
$ws.Range("D13").Value = 1.037689085184316
$ws.Range("E13").Value = 1.033949536831209
$ws.Range("F13").Value = 1.030853219562464
$ws.Range("I13").Value = 1.039293156787416
$ws.Range("J13").Value = 1.041497155657989
$ws.Range("K13").Value = 1.041285206582549
$ws.Range("L13").Value = 1.037559673223235
$ws.Range("M13").Value = 1.034475045518973
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.035093206086204
$ws.Range("D14").Value = 1.037866128033265
$ws.Range("E14").Value = 1.034149818179458
$ws.Range("F14").Value = 1.031195798411578
$ws.Range("I14").Value = 1.039375183025125
$ws.Range("J14").Value = 1.041667927065723
$ws.Range("K14").Value = 1.041430040443885
$ws.Range("L14").Value = 1.037727522702259
$ws.Range("M14").Value = 1.034784544514906
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.035236251565664
$ws.Range("D15").Value = 1.037975169438346
$ws.Range("E15").Value = 1.034273180639171
$ws.Range("F15").Value = 1.031406811537665
$ws.Range("I15").Value = 1.039425665771798
$ws.Range("J15").Value = 1.041773088320245
$ws.Range("K15").Value = 1.041519222274706
$ws.Range("L15").Value = 1.037830889134816
$ws.Range("M15").Value = 1.034975167820905
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.03606812585903
$ws.Range("D16").Value = 1.038609222748328
$ws.Range("E16").Value = 1.034990635862002
$ws.Range("F16").Value = 1.03263408116009
$ws.Range("I16").Value = 1.039718639767267
$ws.Range("J16").Value = 1.042384312224515
$ws.Range("K16").Value = 1.042037461909385
$ws.Range("L16").Value = 1.038431749646674
$ws.Range("M16").Value = 1.03608363289873
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.036589318114723
$ws.Range("D17").Value = 1.039006410655934
$ws.Range("E17").Value = 1.035440181838917
$ws.Range("F17").Value = 1.033403121124179
$ws.Range("I17").Value = 1.039901662282845
$ws.Range("J17").Value = 1.042766964608901
$ws.Range("K17").Value = 1.04236180694764
$ws.Range("L17").Value = 1.038807974295169
$ws.Range("M17").Value = 1.036778036665954
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.03689309978454
$ws.Range("D18").Value = 1.039237891892187
$ws.Range("E18").Value = 1.035702218151954
$ws.Range("F18").Value = 1.033851408800259
$ws.Range("I18").Value = 1.040008146205442
$ws.Range("J18").Value = 1.042989890253042
$ws.Range("K18").Value = 1.042550729469562
$ws.Range("L18").Value = 1.039027176891571
$ws.Range("M18").Value = 1.037182748945907
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.036996644255386
$ws.Range("D19").Value = 1.039316788653007
$ws.Range("E19").Value = 1.035791536059201
$ws.Range("F19").Value = 1.03400421630555
$ws.Range("I19").Value = 1.040044408776438
$ws.Range("J19").Value = 1.043065856700581
$ws.Range("K19").Value = 1.042615102810815
$ws.Range("L19").Value = 1.039101878296381
$ws.Range("M19").Value = 1.037320691299052
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.036533422034409
$ws.Range("D20").Value = 1.038963816027171
$ws.Range("E20").Value = 1.035391968087812
$ws.Range("F20").Value = 1.033320639542132
$ws.Range("I20").Value = 1.039882053653762
$ws.Range("J20").Value = 1.042725937503237
$ws.Range("K20").Value = 1.042327034980832
$ws.Range("L20").Value = 1.038767634117387
$ws.Range("M20").Value = 1.036703567107399
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.035024826149078
$ws.Range("D21").Value = 1.037814001781791
$ws.Range("E21").Value = 1.034090848135659
$ws.Range("F21").Value = 1.031094930221275
$ws.Range("I21").Value = 1.039351040082807
$ws.Range("J21").Value = 1.041617650954883
$ws.Range("K21").Value = 1.041387401968057
$ws.Range("L21").Value = 1.037678105885247
$ws.Range("M21").Value = 1.034693419228413
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.034074264189099
$ws.Range("D22").Value = 1.037089302441054
$ws.Range("E22").Value = 1.033271150675077
$ws.Range("F22").Value = 1.029692889174025
$ws.Range("I22").Value = 1.039014714278975
$ws.Range("J22").Value = 1.040918359145933
$ws.Range("K22").Value = 1.040794216293977
$ws.Range("L22").Value = 1.036990846843737
$ws.Range("M22").Value = 1.033426549442882
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.03457837197207
$ws.Range("D23").Value = 1.037473648474434
$ws.Range("E23").Value = 1.033705844613078
$ws.Range("F23").Value = 1.03043639481504
$ws.Range("I23").Value = 1.039193240965555
$ws.Range("J23").Value = 1.041289303077601
$ws.Range("K23").Value = 1.041108904628359
$ws.Range("L23").Value = 1.037355389200486
$ws.Range("M23").Value = 1.034098431858137
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.036558679749094
$ws.Range("D24").Value = 1.038983063298488
$ws.Range("E24").Value = 1.035413754350674
$ws.Range("F24").Value = 1.033357910282995
$ws.Range("I24").Value = 1.039890914783201
$ws.Range("J24").Value = 1.042744476716957
$ws.Range("K24").Value = 1.042342747746373
$ws.Range("L24").Value = 1.038785862857247
$ws.Range("M24").Value = 1.036737217687711
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.038848077759144
$ws.Range("D25").Value = 1.040727152142897
$ws.Range("E25").Value = 1.03738879653661
$ws.Range("F25").Value = 1.036737232004906
$ws.Range("I25").Value = 1.040689946237862
$ws.Range("J25").Value = 1.044422593942828
$ws.Range("K25").Value = 1.043764276874555
$ws.Range("L25").Value = 1.040436341041456
$ws.Range("M25").Value = 1.039786818744906
